# Fix trading UI (#104)
# The "Recipe" sheet's availableCount column (F) was left at 0 for every
# craftable recipe, which hid them from the in-game trading/crafting UI.
# Flip it to 1 so all recipes show up as available.

$wb = $excel.ActiveWorkbook

$recipe = $wb.Worksheets.Item("Recipe")
$recipe.Range("F2:F21").Value = 1

# Reflect the author's final UI state: they ended the session with the
# "ItemDatas" tab open/scrolled (but no longer the active tab) and the
# "Recipe" tab active with J16 selected.
$itemDatas = $wb.Worksheets.Item("ItemDatas")
$itemDatas.Activate()

$recipe.Activate()
[void]$recipe.Range("J16").Select()
